$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.499.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.914.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2894"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.924.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07554"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.281"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6691"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.485.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007546"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.161.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("E22").Value = "  +5.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.448"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.431"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.70"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.123"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1055"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("E30").Value = "  +2.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.153"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04991"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7295"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9992"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02032"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.016"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.37%  "

$ws.Range("E42").Value = "  +5.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8653"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.822"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.359"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.259"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.15%  "

$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.468"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.68%  "
